$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 19 data: "V17: Distance to pellet includes distance to node"
$ws.Range("A19").Value = "v17-700"
$ws.Range("B19").Value = 2749
$ws.Range("C19").Value = "Distance to pellet includes distance to node"
$ws.Range("D19").Value = 58
$ws.Range("E19").Value = 0.52
$ws.Range("F19").Value = 1.12
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 2588.79
$ws.Range("J19").Value = 1340
$ws.Range("K19").Value = 3380
$ws.Range("L19").Value = 229.02
$ws.Range("M19").Value = 130
$ws.Range("N19").Value = 244
$ws.Range("O19").Value = 0.69
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = 100.22
$ws.Range("S19").Value = 43.3
$ws.Range("T19").Value = 155.9

# Move selection to match the saved workbook state
$ws.Range("E20").Select()
